# Add a new "2nd entry" paragraph (mirroring the existing "1st entry"
# paragraph's structure: a plain run, a superscript run for the ordinal
# suffix, and a plain run) right after the first paragraph, followed by
# a new blank paragraph - matching the target diff.

$d = $word.ActiveDocument

# Step 1: grow the existing lone paragraph's trailing " entry" text into
# " entry" + a new paragraph break + "2nd entry" by using Find/Replace
# (rather than Range.InsertAfter/InsertParagraphAfter) so the freshly
# created run does not inherit stray direct character formatting (this
# host otherwise carries over the last directly-applied run format —
# here, the superscript used for "st" — onto brand new text).
$rngAll = $d.Content
$rngAll.Find.ClearFormatting()
$rngAll.Find.Replacement.ClearFormatting()
$rngAll.Find.Execute("entry", $false, $false, $false, $false, $false, $true, 1, $false, "entry`r2nd entry", 2) | Out-Null

# Step 2: the new paragraph is now Paragraphs(2) ("2nd entry"), all in a
# single plain run. Make just the "nd" ordinal suffix superscript, same
# as "st" in the first paragraph, scoping the Find to that paragraph so
# only this occurrence of "nd" is touched.
$p2 = $d.Paragraphs(2).Range
$p2.Find.ClearFormatting()
$p2.Find.Replacement.ClearFormatting()
$p2.Find.Replacement.Font.Superscript = $true
$p2.Find.Execute("nd", $false, $false, $false, $false, $false, $true, 1, $false, "nd", 2) | Out-Null

# Step 3: append one more, empty, trailing paragraph after "2nd entry"
# (again via Find/Replace, scoped to paragraph 2, to avoid the same
# formatting leak when creating the new paragraph mark).
$p2b = $d.Paragraphs(2).Range
$p2b.Find.ClearFormatting()
$p2b.Find.Replacement.ClearFormatting()
$p2b.Find.Execute(" entry", $false, $false, $false, $false, $false, $true, 1, $false, " entry`r", 2) | Out-Null
